# Update crypto price/volume snapshot (GitHub Actions style refresh).
# Column D ("Price") values are stored as plain text (e.g. "62.083.44",
# "0.999") rather than numbers, so purely-numeric-looking ones are written
# with a leading apostrophe to force Excel to keep them as text instead of
# auto-converting to a number. Column E ("Volume(1h)") values already
# contain padding spaces/percent signs so they stay text without help.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.083.44"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "'3.430.88"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'409.74"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'130.24"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  +6.52%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E9").Value = "  +7.12%  "
$ws.Range("E10").Value = "  +5.06%  "
$ws.Range("D11").Value = "'42.96"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("D12").Value = "'0.0000228"
$ws.Range("E12").Value = "  +52.68%  "
$ws.Range("E13").Value = "  +10.33%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "'3.974.80"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "'21.37"
$ws.Range("E16").Value = "  +7.68%  "
$ws.Range("D17").Value = "'3.439.34"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("E18").Value = "  +7.69%  "
$ws.Range("E19").Value = "  +7.78%  "
$ws.Range("D20").Value = "'62.040.23"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "'456.73"
$ws.Range("E21").Value = "  +46.16%  "
$ws.Range("D22").Value = "'91.66"
$ws.Range("E22").Value = "  +9.28%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  +2.74%  "
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("D26").Value = "'33.18"
$ws.Range("E26").Value = "  +11.77%  "
$ws.Range("D27").Value = "'9.09"
$ws.Range("E27").Value = "  +11.23%  "
$ws.Range("D28").Value = "'4.79"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").Value = "'12.14"
$ws.Range("E30").Value = "  +6.74%  "
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'43.16"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.114"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  +3.37%  "
$ws.Range("D37").Value = "'54.34"
$ws.Range("E37").Value = "  +4.86%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("E40").Value = "  +7.76%  "
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").Value = "'0.319"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "'142.28"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  +9.44%  "
$ws.Range("D45").Value = "'2.01"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").Value = "'2.55"
$ws.Range("E46").Value = "  +15.17%  "
$ws.Range("D47").Value = "'16.68"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").Value = "'22.44"
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("D49").Value = "'0.140"
$ws.Range("E49").Value = "  +17.94%  "
$ws.Range("D50").Value = "'3.776.93"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("E51").Value = "  +8.30%  "
